$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 586, shifting rows 586:685 down to 587:686
$ws.Rows.Item(586).Insert()

# Fill the new row 586 with data
$ws.Cells.Item(586, 1).Value = 4
$ws.Cells.Item(586, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(586, 3).Value = "Los Lagos"
$ws.Cells.Item(586, 4).Value = 45209
$ws.Cells.Item(586, 5).Value = 10
$ws.Cells.Item(586, 6).Value = "Fruta"
$ws.Cells.Item(586, 7).Value = 100102
$ws.Cells.Item(586, 8).Value = "Cítricos"
$ws.Cells.Item(586, 9).Value = 100102006
$ws.Cells.Item(586, 10).Value = "Pomelo"
$ws.Cells.Item(586, 11).Value = "Start Ruby"
$ws.Cells.Item(586, 12).Value = "Primera"
$ws.Cells.Item(586, 13).Value = 100
$ws.Cells.Item(586, 14).Value = 15000
$ws.Cells.Item(586, 15).Value = 15000
$ws.Cells.Item(586, 16).Value = 15000
$ws.Cells.Item(586, 17).Value = "`$/caja 14 kilos empedrada"
$ws.Cells.Item(586, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(586, 19).Value = 1071
$ws.Cells.Item(586, 20).Value = 14

# Apply the same date number format as other date cells in column D
$ws.Cells.Item(586, 4).NumberFormat = $ws.Cells.Item(587, 4).NumberFormat
